$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Route de Darnétal,  Bonsecours"
$ws.Range("C2").Value = 49.4285697
$ws.Range("D2").Value = 1.1453215
$ws.Range("E2").Value = "https://www.google.com/maps/search/?api=1&query=49.4285697,1.1453215"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Quartier du Mont Gargan,  Bonsecours"
$ws.Range("C3").Value = 49.42307500911448
$ws.Range("D3").Value = 1.113333097556789
$ws.Range("E3").Value = "https://www.google.com/maps/search/?api=1&query=49.423075009114484,1.1133330975567892"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Rue Madame de Staël,  Le Châtelet"
$ws.Range("C4").Value = 49.4594655
$ws.Range("D4").Value = 1.1360113
$ws.Range("E4").Value = "https://www.google.com/maps/search/?api=1&query=49.4594655,1.1360113"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Rue Bellevue,  Le Mont Fortin"
$ws.Range("C5").Value = 49.45743925693177
$ws.Range("D5").Value = 1.09952913611965
$ws.Range("E5").Value = "https://www.google.com/maps/search/?api=1&query=49.45743925693177,1.0995291361196504"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Quai Richard Waddington,  Presqu'Île Rollet"
$ws.Range("C6").Value = 49.44372233586422
$ws.Range("D6").Value = 1.05654883905663
$ws.Range("E6").Value = "https://www.google.com/maps/search/?api=1&query=49.44372233586422,1.0565488390566298"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Mosquée El-Fath,  Rue Le Verrier"
$ws.Range("C7").Value = 49.45449635
$ws.Range("D7").Value = 1.141399220863771
$ws.Range("E7").Value = "https://www.google.com/maps/search/?api=1&query=49.45449635,1.1413992208637707"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Stade Irène Hermel,  Allée Jacques Willig"
$ws.Range("C8").Value = 49.42395645
$ws.Range("D8").Value = 1.0987216635178
$ws.Range("E8").Value = "https://www.google.com/maps/search/?api=1&query=49.42395645,1.0987216635177997"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Rue Léonard Bordes,  Quartier du Mont Gargan"
$ws.Range("C9").Value = 49.4343074
$ws.Range("D9").Value = 1.1157783
$ws.Range("E9").Value = "https://www.google.com/maps/search/?api=1&query=49.4343074,1.1157783"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Rue de Fontenelle,  Déville-lès-Rouen"
$ws.Range("C10").Value = 49.47296713563828
$ws.Range("D10").Value = 1.060327903948913
$ws.Range("E10").Value = "https://www.google.com/maps/search/?api=1&query=49.47296713563828,1.0603279039489126"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "La table du Boucher,  3980"
$ws.Range("C11").Value = 49.4747228
$ws.Range("D11").Value = 1.1256959
$ws.Range("E11").Value = "https://www.google.com/maps/search/?api=1&query=49.4747228,1.1256959"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Rue Abbé de l'Épée,  Croix de Pierre"
$ws.Range("C12").Value = 49.4428704
$ws.Range("D12").Value = 1.101857240774189
$ws.Range("E12").Value = "https://www.google.com/maps/search/?api=1&query=49.442870400000004,1.1018572407741893"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Route de Mesnil-Esnard,  Côteaux du Trianon"
$ws.Range("C13").Value = 49.42926201468572
$ws.Range("D13").Value = 1.152774326633283
$ws.Range("E13").Value = "https://www.google.com/maps/search/?api=1&query=49.42926201468572,1.1527743266332833"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Rue Alfred Kastler,  Parc d'activités technologiques La Vatine"
$ws.Range("C14").Value = 49.47288926551011
$ws.Range("D14").Value = 1.096211263258272
$ws.Range("E14").Value = "https://www.google.com/maps/search/?api=1&query=49.47288926551011,1.0962112632582721"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Résidence Muchedent,  Darnétal"
$ws.Range("C15").Value = 49.4501356
$ws.Range("D15").Value = 1.1564396
$ws.Range("E15").Value = "https://www.google.com/maps/search/?api=1&query=49.4501356,1.1564396"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Rue Henri Frère,  Parc de la Varenne"
$ws.Range("C16").Value = 49.46044837309679
$ws.Range("D16").Value = 1.066894627953609
$ws.Range("E16").Value = "https://www.google.com/maps/search/?api=1&query=49.46044837309679,1.0668946279536091"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Rue du Cantony,  Centre Commercial de l'Aubette"
$ws.Range("C17").Value = 49.43451528263126
$ws.Range("D17").Value = 1.155047846947832
$ws.Range("E17").Value = "https://www.google.com/maps/search/?api=1&query=49.43451528263126,1.1550478469478322"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Pont Flaubert,  Quai de France"
$ws.Range("C18").Value = 49.44133069999999
$ws.Range("D18").Value = 1.065146643849046
$ws.Range("E18").Value = "https://www.google.com/maps/search/?api=1&query=49.441330699999995,1.0651466438490456"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Chemin de la Grand'Mare,  Vallon Suisse"
$ws.Range("C19").Value = 49.45606676013594
$ws.Range("D19").Value = 1.134587443855071
$ws.Range("E19").Value = "https://www.google.com/maps/search/?api=1&query=49.45606676013594,1.134587443855071"

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Boulevard Industriel,  Sotteville-lès-Rouen"
$ws.Range("C20").Value = 49.40399859811959
$ws.Range("D20").Value = 1.110365343518302
$ws.Range("E20").Value = "https://www.google.com/maps/search/?api=1&query=49.40399859811959,1.110365343518302"

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Rue de l'Appel du 18 Juin 1940,  La Prévotière"
$ws.Range("C21").Value = 49.46639417381033
$ws.Range("D21").Value = 1.12582431140331
$ws.Range("E21").Value = "https://www.google.com/maps/search/?api=1&query=49.46639417381033,1.1258243114033095"

$ws.Rows.Item(22).Delete()
